$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(32, 8).Value = 4277.75
$ws.Cells.Item(32, 9).Value = 1609.5
$ws.Cells.Item(32, 10).Value = 5167.1665
$ws.Cells.Item(32, 11).Value = 1609.5
$ws.Cells.Item(32, 12).Value = 5167.1665
$ws.Cells.Item(32, 13).Value = -1283.5
$ws.Cells.Item(32, 14).Value = -5819.1665

$ws.Cells.Item(86, 8).Value = 144630.86
$ws.Cells.Item(86, 9).Value = 1716.125
$ws.Cells.Item(86, 11).Value = 1716.125
$ws.Cells.Item(86, 13).Value = -593.125

$ws.Cells.Item(88, 8).Value = 1583.9286
$ws.Cells.Item(88, 9).Value = 1707
$ws.Cells.Item(88, 10).Value = 1491.625
$ws.Cells.Item(88, 11).Value = 1707
$ws.Cells.Item(88, 12).Value = 1491.625
$ws.Cells.Item(88, 13).Value = -1301
$ws.Cells.Item(88, 14).Value = -2303.625

$ws.Cells.Item(89, 8).Value = 144630.86
$ws.Cells.Item(89, 9).Value = 1716.125
$ws.Cells.Item(89, 11).Value = 8580.625
$ws.Cells.Item(89, 13).Value = -2964.625

$ws.Cells.Item(91, 8).Value = 1583.9286
$ws.Cells.Item(91, 9).Value = 1707
$ws.Cells.Item(91, 10).Value = 1491.625
$ws.Cells.Item(91, 11).Value = 1707
$ws.Cells.Item(91, 12).Value = 1491.625
$ws.Cells.Item(91, 13).Value = -303
$ws.Cells.Item(91, 14).Value = -4299.625

$ws.Cells.Item(111, 8).Value = 13909.786
$ws.Cells.Item(111, 10).Value = 8688
$ws.Cells.Item(111, 12).Value = 26064
$ws.Cells.Item(111, 14).Value = -32198

$ws.Cells.Item(132, 8).Value = 5400
$ws.Cells.Item(132, 9).Value = 5400
$ws.Cells.Item(132, 11).Value = 16200
$ws.Cells.Item(132, 13).Value = -13670

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(23, 8).Value = 39999
$ws.Cells.Item(23, 10).Value = 39999
$ws.Cells.Item(23, 12).Value = 39999
$ws.Cells.Item(23, 14).Value = -40517

$ws.Cells.Item(32, 8).Value = 8561919
$ws.Cells.Item(32, 9).Value = 8757158
$ws.Cells.Item(32, 11).Value = 8757158
$ws.Cells.Item(32, 13).Value = -8756871

$ws.Cells.Item(88, 8).Value = 1260.55
$ws.Cells.Item(88, 9).Value = 1158.3334
$ws.Cells.Item(88, 10).Value = 1304.3572
$ws.Cells.Item(88, 11).Value = 1158.3334
$ws.Cells.Item(88, 12).Value = 1304.3572
$ws.Cells.Item(88, 13).Value = -752.3334
$ws.Cells.Item(88, 14).Value = -2116.3572

$ws.Cells.Item(91, 8).Value = 1260.55
$ws.Cells.Item(91, 9).Value = 1158.3334
$ws.Cells.Item(91, 10).Value = 1304.3572
$ws.Cells.Item(91, 11).Value = 1158.3334
$ws.Cells.Item(91, 12).Value = 1304.3572
$ws.Cells.Item(91, 13).Value = 245.6666
$ws.Cells.Item(91, 14).Value = -4112.3572

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 1144.4615
$ws.Cells.Item(86, 9).Value = 934.5
$ws.Cells.Item(86, 10).Value = 1480.4
$ws.Cells.Item(86, 11).Value = 934.5
$ws.Cells.Item(86, 12).Value = 1480.4
$ws.Cells.Item(86, 13).Value = 188.5
$ws.Cells.Item(86, 14).Value = -3726.4

$ws.Cells.Item(89, 8).Value = 1144.4615
$ws.Cells.Item(89, 9).Value = 934.5
$ws.Cells.Item(89, 10).Value = 1480.4
$ws.Cells.Item(89, 11).Value = 4672.5
$ws.Cells.Item(89, 12).Value = 7402
$ws.Cells.Item(89, 13).Value = 943.5
$ws.Cells.Item(89, 14).Value = -18634

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(19, 8).Value = 90.5
$ws.Cells.Item(19, 9).Value = 90.5
$ws.Cells.Item(19, 10).Value = 0
$ws.Cells.Item(19, 11).Value = 90.5
$ws.Cells.Item(19, 12).Value = 0
$ws.Cells.Item(19, 13).Value = 79.5
$ws.Cells.Item(19, 14).Value = $null

$ws.Cells.Item(24, 8).Value = 90.5
$ws.Cells.Item(24, 9).Value = 90.5
$ws.Cells.Item(24, 10).Value = 0
$ws.Cells.Item(24, 11).Value = 90.5
$ws.Cells.Item(24, 12).Value = 0
$ws.Cells.Item(24, 13).Value = 79.5
$ws.Cells.Item(24, 14).Value = $null

$ws.Cells.Item(31, 8).Value = 2548.75
$ws.Cells.Item(31, 10).Value = 4222
$ws.Cells.Item(31, 12).Value = 4222
$ws.Cells.Item(31, 14).Value = -4812

$ws.Cells.Item(34, 8).Value = 2548.75
$ws.Cells.Item(34, 10).Value = 4222
$ws.Cells.Item(34, 12).Value = 4222
$ws.Cells.Item(34, 14).Value = -4626

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(120, 8).Value = 1999
$ws.Cells.Item(120, 9).Value = 1999
$ws.Cells.Item(120, 11).Value = 5997
$ws.Cells.Item(120, 13).Value = -1159

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(63, 8).Value = 38999
$ws.Cells.Item(63, 10).Value = 38999
$ws.Cells.Item(63, 12).Value = 38999
$ws.Cells.Item(63, 14).Value = -40371

$ws.Cells.Item(66, 8).Value = 38999
$ws.Cells.Item(66, 10).Value = 38999
$ws.Cells.Item(66, 12).Value = 116997
$ws.Cells.Item(66, 14).Value = -123861

$ws.Cells.Item(80, 8).Value = 1574.75
$ws.Cells.Item(80, 9).Value = 3000
$ws.Cells.Item(80, 11).Value = 3000
$ws.Cells.Item(80, 13).Value = -2002

$ws.Cells.Item(83, 8).Value = 1574.75
$ws.Cells.Item(83, 9).Value = 3000
$ws.Cells.Item(83, 11).Value = 15000
$ws.Cells.Item(83, 13).Value = -10008

$ws.Cells.Item(102, 8).Value = 1898.6
$ws.Cells.Item(102, 9).Value = 1898.6
$ws.Cells.Item(102, 10).Value = 0
$ws.Cells.Item(102, 11).Value = 1898.6
$ws.Cells.Item(102, 12).Value = 0
$ws.Cells.Item(102, 13).Value = -276.5999999999999
$ws.Cells.Item(102, 14).Value = $null

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 1998.5
$ws.Cells.Item(22, 9).Value = 0
$ws.Cells.Item(22, 10).Value = 1998.5
$ws.Cells.Item(22, 11).Value = 0
$ws.Cells.Item(22, 12).Value = 1998.5
$ws.Cells.Item(22, 13).Value = $null
$ws.Cells.Item(22, 14).Value = -2588.5

$ws.Cells.Item(27, 8).Value = 1998.5
$ws.Cells.Item(27, 9).Value = 0
$ws.Cells.Item(27, 10).Value = 1998.5
$ws.Cells.Item(27, 11).Value = 0
$ws.Cells.Item(27, 12).Value = 1998.5
$ws.Cells.Item(27, 13).Value = $null
$ws.Cells.Item(27, 14).Value = -2212.5

$ws.Cells.Item(82, 8).Value = 2791.3333
$ws.Cells.Item(82, 10).Value = 3499
$ws.Cells.Item(82, 12).Value = 3499
$ws.Cells.Item(82, 14).Value = -4221

$ws.Cells.Item(85, 8).Value = 2791.3333
$ws.Cells.Item(85, 10).Value = 3499
$ws.Cells.Item(85, 12).Value = 3499
$ws.Cells.Item(85, 14).Value = -5995

$ws.Cells.Item(100, 8).Value = 3359.3
$ws.Cells.Item(100, 9).Value = 3121.4443
$ws.Cells.Item(100, 11).Value = 3121.4443
$ws.Cells.Item(100, 13).Value = -2580.4443

$ws.Cells.Item(136, 8).Value = 2281.5
$ws.Cells.Item(136, 9).Value = 2299.8
$ws.Cells.Item(136, 11).Value = 6899.400000000001
$ws.Cells.Item(136, 13).Value = -4349.400000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(41, 8).Value = 14038.833
$ws.Cells.Item(41, 9).Value = 10139.667
$ws.Cells.Item(41, 11).Value = 10139.667
$ws.Cells.Item(41, 13).Value = -9749.666999999999

$ws.Cells.Item(136, 8).Value = 2208.1667
$ws.Cells.Item(136, 9).Value = 1566.3334
$ws.Cells.Item(136, 11).Value = 4699.0002
$ws.Cells.Item(136, 13).Value = -2149.0002
